$d = $word.ActiveDocument

# ---- Paragraph 1: Title ----
$p1 = $d.Paragraphs(1)
$start1 = $p1.Range.Start
$end1 = $p1.Range.End
$words1 = @("Questions:", " ", "Introduction", " ", "to", " ", "simultaneous", " ", "equations")
$runsXml1 = ""
foreach ($w in $words1) {
    $runsXml1 += "<w:r><w:t xml:space=`"preserve`">$w</w:t></w:r>"
}
$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$runsXml1</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$target1 = $d.Range($start1, $end1 - 1)
$target1.InsertXML($xml1)

# ---- Paragraph 2: Author ----
$p2 = $d.Paragraphs(2)
$start2 = $p2.Range.Start
$end2 = $p2.Range.End
$words2 = @("Ollie", " ", "Brooke")
$runsXml2 = ""
foreach ($w in $words2) {
    $runsXml2 += "<w:r><w:t xml:space=`"preserve`">$w</w:t></w:r>"
}
$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$runsXml2</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$target2 = $d.Range($start2, $end2 - 1)
$target2.InsertXML($xml2)

# ---- Paragraph 4: Abstract ----
$p4 = $d.Paragraphs(4)
$start4 = $p4.Range.Start
$end4 = $p4.Range.End
$words4 = @("Questions", " ", "relating", " ", "to", " ", "the", " ", "introduction", " ", "to", " ", "simultaneous", " ", "equations", " ", "study", " ", "guide.")
$runsXml4 = ""
foreach ($w in $words4) {
    $runsXml4 += "<w:r><w:t xml:space=`"preserve`">$w</w:t></w:r>"
}
$xml4 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$runsXml4</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$target4 = $d.Range($start4, $end4 - 1)
$target4.InsertXML($xml4)
